# Update database and shift read_price algorithm by one year:
# drop the oldest (1396/12) period column, shift remaining periods/dates/
# financial figures one column to the left, and populate the newly freed
# rightmost column (H) with the new 1401/12 period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: financial period headers ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates ---
$ws.Range("D9").Value = "1399-02-28 (9)"
$ws.Range("E9").Value = "1400-04-16 (8)"
$ws.Range("F9").Value = "1401-03-24 (9)"
$ws.Range("G9").Value = "1402-02-27 (7)"
# H9 is a bare date-looking string ("1402-02-27"); a plain assignment gets
# auto-parsed into a date serial number by the "smart" value setter, so
# force it in as literal text with a leading quote, then restore the
# cell's original formatting (the quote-prefix entry tweaks the style).
$ws.Range("H9").Value = "'1402-02-27"
$ws.Range("G9").Copy()
$ws.Range("H9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 11: فروش (Sales) ---
$ws.Range("D11").Value = 7976
$ws.Range("E11").Value = 8995
$ws.Range("F11").Value = 4459
$ws.Range("G11").Value = 6186
$ws.Range("H11").Value = 8500

# --- Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold) ---
$ws.Range("D12").Value = -6468
$ws.Range("E12").Value = -6847
$ws.Range("F12").Value = -3306
$ws.Range("G12").Value = -4319
$ws.Range("H12").Value = -6766

# --- Row 13: سود (زیان) ناخالص (Gross profit) ---
$ws.Range("D13").Value = 1508
$ws.Range("E13").Value = 2148
$ws.Range("F13").Value = 1153
$ws.Range("G13").Value = 1867
$ws.Range("H13").Value = 1734

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) ---
$ws.Range("D14").Value = -49
$ws.Range("E14").Value = -150
$ws.Range("F14").Value = -146
$ws.Range("G14").Value = -230
$ws.Range("H14").Value = -270

# Row 15 (impairment expense) is "-" across all periods in both the
# before and after states, so no write is required there.

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other op. income/exp) ---
$ws.Range("D16").Value = 9
$ws.Range("E16").Value = 19
$ws.Range("F16").Value = 19
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 97

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ---
$ws.Range("D17").Value = 1468
$ws.Range("E17").Value = 2017
$ws.Range("F17").Value = 1026
$ws.Range("G17").Value = 1652
$ws.Range("H17").Value = 1561

# --- Row 18: هزینه های مالی (Financial expenses) ---
$ws.Range("D18").Value = -156
$ws.Range("E18").Value = -118
$ws.Range("F18").Value = -197
$ws.Range("G18").Value = -640
$ws.Range("H18").Value = -745

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-op. income/exp) ---
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 6
$ws.Range("F19").Value = 7
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = 83

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit) ---
$ws.Range("D20").Value = 1319
$ws.Range("E20").Value = 1905
$ws.Range("F20").Value = 835
$ws.Range("G20").Value = 1069
$ws.Range("H20").Value = 899

# --- Row 21: مالیات (Tax) ---
$ws.Range("D21").Value = -219
$ws.Range("E21").Value = -261
$ws.Range("F21").Value = -145
$ws.Range("G21").Value = -216
$ws.Range("H21").Value = -113

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم (Net continuing profit) ---
$ws.Range("D22").Value = 1101
$ws.Range("E22").Value = 1644
$ws.Range("F22").Value = 690
$ws.Range("G22").Value = 853
$ws.Range("H22").Value = 785

# Row 23 (discontinued operations) is "-" across all periods in both the
# before and after states, so no write is required there.

# --- Row 24: سود (زیان) خالص (Net profit) ---
$ws.Range("D24").Value = 1101
$ws.Range("E24").Value = 1644
$ws.Range("F24").Value = 690
$ws.Range("G24").Value = 853
$ws.Range("H24").Value = 785

# Row 25 (EPS after tax) is 0 across all periods in both the before and
# after states, so no write is required there.

# --- Row 26: سرمایه (Capital) ---
$ws.Range("D26").Value = 989
$ws.Range("E26").Value = 780
$ws.Range("F26").Value = 3428
$ws.Range("G26").Value = 2937
$ws.Range("H26").Value = 2196

# Row 27 (EPS on latest capital) is 0 across all periods in both the
# before and after states, so no write is required there.
